# Append 4 new weekly price rows (392-395) for "Pepino dulce" at the
# Mercado Mayorista Lo Valledor de Santiago, following the same layout
# as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerial = 45075
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$newRows = @(
    @{ Row = 392; I = "Especial"; J = 300; K = 14000; L = 14000; M = 14000; P = 778 },
    @{ Row = 393; I = "Primera";  J = 480; K = 12000; L = 12000; M = 12000; P = 667 },
    @{ Row = 394; I = "Segunda";  J = 370; K = 10000; L = 10000; M = 10000; P = 556 },
    @{ Row = 395; I = "Tercera";  J = 220; K = 8000;  L = 8000;  M = 8000;  P = 444 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"

    $ws.Cells.Item($row, 4).Value = $dateSerial
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
